$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7 with Merge Sort details
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Merge Sort"
$ws.Range("D7").Value = "c"
$ws.Range("E7").Value = "O(nlogn)"
$ws.Range("F7").Value = "O(n)"
$ws.Range("H7").Value = "MergeSort ( Stable sorting algo)"

# Move the active selection to A8, mirroring the post-edit cursor position
$ws.Range("A8").Select()
